# 4.b.1.1 worksheet update: add a new "2020" data column (J) to the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (thin separator row above the header): new bottom-bordered blank cell J3.
$ws.Range("J3").Borders.Item(9).Weight = -4138

# Row 4 (year header row): J4 = 2020, formatted like the existing year cells.
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = 2020

# Row 5: J5 = 370
$ws.Range("I5").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("J5").Value = 370

# Row 6: J6 = 5
$ws.Range("I6").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("J6").Value = 5

# Row 7: J7 = 5
$ws.Range("I7").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J7").Value = 5

# Row 8: I8 was the placeholder "-" and is now a real value (42); J8 = 20
$ws.Range("I8").Value = 42
$ws.Range("I8").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("J8").Value = 20

# Row 9: I9 365->30 (updated total), J9 = 19
$ws.Range("I9").Value = 30
$ws.Range("I9").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("J9").Value = 19

# Row 10 (totals row): I10 67->62 (updated total), J10 = 73
$ws.Range("I10").Value = 62
$ws.Range("I10").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("J10").Value = 73
